$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 52.13946256887639
$ws.Range("D2").Value = 0.4971469653530965

$ws.Range("C3").Value = 163107388526.8221
$ws.Range("D3").Value = -1573070400.055565

$ws.Range("C4").Value = 53.89935221717534
$ws.Range("D4").Value = 0.4801739125694827

$ws.Range("A5").Value = "Lasso Regression with positive Coefficients"
$ws.Range("C5").Value = 73.95363565634419
$ws.Range("D5").Value = 0.2867626883603009

$ws.Range("C6").Value = 23.98751883576024
$ws.Range("D6").Value = 0.7686551405420108

$ws.Range("C7").Value = 14.2444248244857
$ws.Range("D7").Value = 0.8626212872767882

$ws.Range("C8").Value = 42.85226970232395
$ws.Range("D8").Value = 0.5867162260666638

$ws.Range("C9").Value = 91.03617693129489
$ws.Range("D9").Value = 0.1875999591520416

$ws.Range("C10").Value = 2749121268259.972
$ws.Range("D10").Value = -24532952786.7245

$ws.Range("C11").Value = 91.95489924413481
$ws.Range("D11").Value = 0.1794013498778146

$ws.Range("A12").Value = "Lasso Regression with positive Coefficients"
$ws.Range("C12").Value = 99.68091203694125
$ws.Range("D12").Value = 0.1104549889909251

$ws.Range("C13").Value = 92.89837300804287
$ws.Range("D13").Value = 0.1709818605036471

$ws.Range("C14").Value = 47.29551685833902
$ws.Range("D14").Value = 0.5779383414064311

$ws.Range("C15").Value = 75.67969437911121
$ws.Range("D15").Value = 0.3246400620343314

$ws.Range("C16").Value = 69.48906955627709
$ws.Range("D16").Value = 0.4187810453552685

$ws.Range("C17").Value = 465674854319.1072
$ws.Range("D17").Value = -3894987423.06517

$ws.Range("C18").Value = 70.52409416323501
$ws.Range("D18").Value = 0.4101239151917907

$ws.Range("A19").Value = "Lasso Regression with positive Coefficients"
$ws.Range("C19").Value = 108.0457870943962
$ws.Range("D19").Value = 0.09628579243079727

$ws.Range("C20").Value = 47.54066055924312
$ws.Range("D20").Value = 0.6023614474937622

$ws.Range("C21").Value = 29.49500624111414
$ws.Range("D21").Value = 0.7532985143682698

$ws.Range("C22").Value = 55.80950074928899
$ws.Range("D22").Value = 0.5331993953599288

$ws.Range("C23").Value = 29.36788088280915
$ws.Range("D23").Value = 0.5250734874983307

$ws.Range("C24").Value = 2319302743.494648
$ws.Range("D24").Value = -37506905.53503171

$ws.Range("C25").Value = 30.27669265763229
$ws.Range("D25").Value = 0.5103765194583296

$ws.Range("A26").Value = "Lasso Regression with positive Coefficients"
$ws.Range("C26").Value = 53.74202139753132
$ws.Range("D26").Value = 0.130903898072532

$ws.Range("C27").Value = 27.62767539645994
$ws.Range("D27").Value = 0.5532154472797014

$ws.Range("C28").Value = 11.60238164041661
$ws.Range("D28").Value = 0.8123705734443384

$ws.Range("C29").Value = 27.95411761370189
$ws.Range("D29").Value = 0.5479363444262577

$ws.Range("C30").Value = 38.43626758216009
$ws.Range("D30").Value = 0.4382436275045587

$ws.Range("C31").Value = 91993686026.33995
$ws.Range("D31").Value = -1344512424.515199

$ws.Range("C32").Value = 38.88037443180772
$ws.Range("D32").Value = 0.4317528866352704

$ws.Range("A33").Value = "Lasso Regression with positive Coefficients"
$ws.Range("C33").Value = 67.0402920782105
$ws.Range("D33").Value = 0.02018812809050807

$ws.Range("C34").Value = 28.60817844625111
$ws.Range("D34").Value = 0.5818837894882529

$ws.Range("C35").Value = 11.6772285141512
$ws.Range("D35").Value = 0.8293341694302656

$ws.Range("C36").Value = 46.30263178075999
$ws.Range("D36").Value = 0.3232746022867063
